$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 101
$ws.Range('A101').Value = '''September25  18:01:38'
$ws.Range('B101').Value = '''resnet18'
$ws.Range('C101').Value = '''1000'
$ws.Range('D101').Value = '''1'
$ws.Range('E101').Value = '''True'
$ws.Range('F101').Value = '''sgd'
$ws.Range('G101').Value = '''0.0xsingle + 1.0Xmulti'
$ws.Range('H101').Value = '''0.01'
$ws.Range('I101').Value = '''0.9'
$ws.Range('J101').Value = '''<function exp_lr_scheduler at 0x7fb7a2408840>'
$ws.Range('K101').Value = '''10'
$ws.Range('L101').Value = '''True'
$ws.Range('M101').Value = '''16'
$ws.Range('S101').Value = ''' '

# Row 102
$ws.Range('A102').Value = '''September25  18:02:04'
$ws.Range('B102').Value = '''resnet18'
$ws.Range('C102').Value = '''1000'
$ws.Range('D102').Value = '''1'
$ws.Range('E102').Value = '''True'
$ws.Range('F102').Value = '''sgd'
$ws.Range('G102').Value = '''0.0xsingle + 1.0Xmulti'
$ws.Range('H102').Value = '''0.01'
$ws.Range('I102').Value = '''0.9'
$ws.Range('J102').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K102').Value = '''10'
$ws.Range('L102').Value = '''True'
$ws.Range('M102').Value = '''16'
$ws.Range('S102').Value = ''' '

# Row 103
$ws.Range('A103').Value = '''September25  18:03:52'
$ws.Range('B103').Value = '''resnet18'
$ws.Range('C103').Value = '''1000'
$ws.Range('D103').Value = '''1'
$ws.Range('E103').Value = '''True'
$ws.Range('F103').Value = '''sgd'
$ws.Range('G103').Value = '''0.0xsingle + 1.0Xmulti'
$ws.Range('H103').Value = '''0.01'
$ws.Range('I103').Value = '''0.9'
$ws.Range('J103').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K103').Value = '''10'
$ws.Range('L103').Value = '''True'
$ws.Range('M103').Value = '''16'
$ws.Range('N103').Value = 13
$ws.Range('O103').Value = 0.0079049546426783
$ws.Range('P103').Value = 0.02518401054897569
$ws.Range('Q103').Value = 0.5755208333333334
$ws.Range('R103').Value = 0.3829787234042553
$ws.Range('S103').Value = 0.733321496116585
$ws.Range('T103').Value = 1.26443038342523
$ws.Range('U103').Value = 13
$ws.Range('V103').Value = 0.5755208333333334
$ws.Range('W103').Value = 0.3829787234042553
$ws.Range('X103').Value = 0.733321496116585
$ws.Range('Y103').Value = 1.26443038342523

# Row 104
$ws.Range('A104').Value = '''September25  18:08:53'
$ws.Range('B104').Value = '''resnet18'
$ws.Range('C104').Value = '''1000'
$ws.Range('D104').Value = '''1'
$ws.Range('E104').Value = '''True'
$ws.Range('F104').Value = '''sgd'
$ws.Range('G104').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H104').Value = '''0.01'
$ws.Range('I104').Value = '''0.9'
$ws.Range('J104').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K104').Value = '''10'
$ws.Range('L104').Value = '''True'
$ws.Range('M104').Value = '''16'
$ws.Range('N104').Value = 20
$ws.Range('O104').Value = 0.04323956785568347
$ws.Range('P104').Value = 0.1327234439574476
$ws.Range('Q104').Value = 0.7421875
$ws.Range('R104').Value = 0.425531914893617
$ws.Range('S104').Value = 0.9319278852643768
$ws.Range('T104').Value = 1.439773707312186
$ws.Range('U104').Value = 20
$ws.Range('V104').Value = 0.7421875
$ws.Range('W104').Value = 0.425531914893617
$ws.Range('X104').Value = 0.9319278852643768
$ws.Range('Y104').Value = 1.439773707312186

# Row 105
$ws.Range('A105').Value = '''September25  18:15:52'
$ws.Range('B105').Value = '''resnet18'
$ws.Range('C105').Value = '''1000'
$ws.Range('D105').Value = '''1'
$ws.Range('E105').Value = '''True'
$ws.Range('F105').Value = '''sgd'
$ws.Range('G105').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H105').Value = '''0.01'
$ws.Range('I105').Value = '''0.9'
$ws.Range('J105').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K105').Value = '''10'
$ws.Range('L105').Value = '''True'
$ws.Range('M105').Value = '''16'
$ws.Range('N105').Value = 4
$ws.Range('O105').Value = 0
$ws.Range('P105').Value = 0
$ws.Range('Q105').Value = 0.13671875
$ws.Range('R105').Value = 0.121580547112462
$ws.Range('S105').Value = 3.552801124934146
$ws.Range('T105').Value = 3.536608373031219
$ws.Range('U105').Value = 4
$ws.Range('V105').Value = 0.13671875
$ws.Range('W105').Value = 0.121580547112462
$ws.Range('X105').Value = 3.552801124934146
$ws.Range('Y105').Value = 3.536608373031219

# Row 106
$ws.Range('A106').Value = '''September25  18:17:25'
$ws.Range('B106').Value = '''resnet18'
$ws.Range('C106').Value = '''1000'
$ws.Range('D106').Value = '''1'
$ws.Range('E106').Value = '''True'
$ws.Range('F106').Value = '''sgd'
$ws.Range('G106').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H106').Value = '''0.01'
$ws.Range('I106').Value = '''0.9'
$ws.Range('J106').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K106').Value = '''10'
$ws.Range('L106').Value = '''True'
$ws.Range('M106').Value = '''16'
$ws.Range('N106').Value = 2
$ws.Range('O106').Value = 0
$ws.Range('P106').Value = 0
$ws.Range('Q106').Value = 0.1041666666666667
$ws.Range('R106').Value = 0.1155015197568389
$ws.Range('S106').Value = 3.219609513072458
$ws.Range('T106').Value = 3.356253519919138
$ws.Range('U106').Value = 2
$ws.Range('V106').Value = 0.1041666666666667
$ws.Range('W106').Value = 0.1155015197568389
$ws.Range('X106').Value = 3.219609513072458
$ws.Range('Y106').Value = 3.356253519919138

# Row 107
$ws.Range('A107').Value = '''September25  18:18:11'
$ws.Range('B107').Value = '''resnet18'
$ws.Range('C107').Value = '''1000'
$ws.Range('D107').Value = '''1'
$ws.Range('E107').Value = '''True'
$ws.Range('F107').Value = '''sgd'
$ws.Range('G107').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H107').Value = '''0.01'
$ws.Range('I107').Value = '''0.9'
$ws.Range('J107').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K107').Value = '''10'
$ws.Range('L107').Value = '''True'
$ws.Range('M107').Value = '''16'
$ws.Range('N107').Value = 12
$ws.Range('O107').Value = 0.002769226039769516
$ws.Range('P107').Value = 0.01200845931574566
$ws.Range('Q107').Value = 0.89453125
$ws.Range('R107').Value = 0.4620060790273556
$ws.Range('S107').Value = 0.625
$ws.Range('T107').Value = 1.241383677454499
$ws.Range('U107').Value = 12
$ws.Range('V107').Value = 0.89453125
$ws.Range('W107').Value = 0.4620060790273556
$ws.Range('X107').Value = 0.625
$ws.Range('Y107').Value = 1.241383677454499

# Row 108
$ws.Range('A108').Value = '''September25  18:22:37'
$ws.Range('B108').Value = '''resnet18'
$ws.Range('C108').Value = '''1000'
$ws.Range('D108').Value = '''1'
$ws.Range('E108').Value = '''True'
$ws.Range('F108').Value = '''sgd'
$ws.Range('G108').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H108').Value = '''0.01'
$ws.Range('I108').Value = '''0.9'
$ws.Range('J108').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K108').Value = '''10'
$ws.Range('L108').Value = '''True'
$ws.Range('M108').Value = '''16'
$ws.Range('N108').Value = 10
$ws.Range('O108').Value = 0.002782941594584069
$ws.Range('P108').Value = 0.004269491130978684
$ws.Range('Q108').Value = 0.4479166666666667
$ws.Range('R108').Value = 0.3860182370820669
$ws.Range('S108').Value = 1.001951221367587
$ws.Range('T108').Value = 1.214151610688962
$ws.Range('U108').Value = 10
$ws.Range('V108').Value = 0.4479166666666667
$ws.Range('W108').Value = 0.3860182370820669
$ws.Range('X108').Value = 1.001951221367587
$ws.Range('Y108').Value = 1.214151610688962

# Row 109
$ws.Range('A109').Value = '''September25  18:27:20'
$ws.Range('B109').Value = '''resnet18'
$ws.Range('C109').Value = '''1000'
$ws.Range('D109').Value = '''1'
$ws.Range('E109').Value = '''True'
$ws.Range('F109').Value = '''sgd'
$ws.Range('G109').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H109').Value = '''0.01'
$ws.Range('I109').Value = '''0.9'
$ws.Range('J109').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K109').Value = '''10'
$ws.Range('L109').Value = '''True'
$ws.Range('M109').Value = '''16'
$ws.Range('N109').Value = 18
$ws.Range('O109').Value = 0.002259189934799603
$ws.Range('P109').Value = 0.004330871043596586
$ws.Range('Q109').Value = 0.6875
$ws.Range('R109').Value = 0.4316109422492401
$ws.Range('S109').Value = 0.7525996611745185
$ws.Range('T109').Value = 1.215402668837128
$ws.Range('U109').Value = 18
$ws.Range('V109').Value = 0.6875
$ws.Range('W109').Value = 0.4316109422492401
$ws.Range('X109').Value = 0.7525996611745185
$ws.Range('Y109').Value = 1.215402668837128

# Row 110
$ws.Range('A110').Value = '''September25  18:32:12'
$ws.Range('B110').Value = '''resnet18'
$ws.Range('C110').Value = '''1000'
$ws.Range('D110').Value = '''1'
$ws.Range('E110').Value = '''True'
$ws.Range('F110').Value = '''sgd'
$ws.Range('G110').Value = '''0.0xsingle + 1.0Xmulti'
$ws.Range('H110').Value = '''0.01'
$ws.Range('I110').Value = '''0.9'
$ws.Range('J110').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K110').Value = '''10'
$ws.Range('L110').Value = '''True'
$ws.Range('M110').Value = '''16'
$ws.Range('N110').Value = 5
$ws.Range('O110').Value = 0.01706836840215449
$ws.Range('P110').Value = 0.02383613747213387
$ws.Range('Q110').Value = 0.4166666666666667
$ws.Range('R110').Value = 0.3586626139817629
$ws.Range('S110').Value = 1.012937148428602
$ws.Range('T110').Value = 1.372770366901603
$ws.Range('U110').Value = 5
$ws.Range('V110').Value = 0.4166666666666667
$ws.Range('W110').Value = 0.3586626139817629
$ws.Range('X110').Value = 1.012937148428602
$ws.Range('Y110').Value = 1.372770366901603

# Row 111
$ws.Range('A111').Value = '''September25  18:33:18'
$ws.Range('B111').Value = '''resnet18'
$ws.Range('C111').Value = '''1000'
$ws.Range('D111').Value = '''1'
$ws.Range('E111').Value = '''True'
$ws.Range('F111').Value = '''sgd'
$ws.Range('G111').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H111').Value = '''0.01'
$ws.Range('I111').Value = '''0.9'
$ws.Range('J111').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K111').Value = '''10'
$ws.Range('L111').Value = '''True'
$ws.Range('M111').Value = '''16'
$ws.Range('N111').Value = 1
$ws.Range('O111').Value = 0.007602096787498643
$ws.Range('P111').Value = 0.006338213635523631
$ws.Range('Q111').Value = 0.2057291666666667
$ws.Range('R111').Value = 0.2644376899696049
$ws.Range('S111').Value = 2.590547722007838
$ws.Range('T111').Value = 2.049538460906641
$ws.Range('U111').Value = 1
$ws.Range('V111').Value = 0.2057291666666667
$ws.Range('W111').Value = 0.2644376899696049
$ws.Range('X111').Value = 2.590547722007838
$ws.Range('Y111').Value = 2.049538460906641

# Row 112
$ws.Range('A112').Value = '''September25  18:33:39'
$ws.Range('B112').Value = '''resnet18'
$ws.Range('C112').Value = '''1000'
$ws.Range('D112').Value = '''1'
$ws.Range('E112').Value = '''True'
$ws.Range('F112').Value = '''sgd'
$ws.Range('G112').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H112').Value = '''0.01'
$ws.Range('I112').Value = '''0.9'
$ws.Range('J112').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K112').Value = '''10'
$ws.Range('L112').Value = '''True'
$ws.Range('M112').Value = '''16'
$ws.Range('S112').Value = ''' '

# Row 113
$ws.Range('A113').Value = '''September25  18:34:54'
$ws.Range('B113').Value = '''resnet18'
$ws.Range('C113').Value = '''1000'
$ws.Range('D113').Value = '''1'
$ws.Range('E113').Value = '''True'
$ws.Range('F113').Value = '''sgd'
$ws.Range('G113').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H113').Value = '''0.01'
$ws.Range('I113').Value = '''0.9'
$ws.Range('J113').Value = '''<function exp_lr_scheduler at 0x7fb7a24087b8>'
$ws.Range('K113').Value = '''10'
$ws.Range('L113').Value = '''True'
$ws.Range('M113').Value = '''16'
$ws.Range('N113').Value = 15
$ws.Range('O113').Value = 0.002433132079507535
$ws.Range('P113').Value = 0.003921142020644932
$ws.Range('Q113').Value = 0.5442708333333334
$ws.Range('R113').Value = 0.4012158054711246
$ws.Range('S113').Value = 0.8853553900364908
$ws.Range('T113').Value = 1.161698787529852
$ws.Range('U113').Value = 15
$ws.Range('V113').Value = 0.5442708333333334
$ws.Range('W113').Value = 0.4012158054711246
$ws.Range('X113').Value = 0.8853553900364908
$ws.Range('Y113').Value = 1.161698787529852

# Row 114
$ws.Range('A114').Value = '''September25  18:42:46'
$ws.Range('B114').Value = '''resnet34'
$ws.Range('C114').Value = '''1000'
$ws.Range('D114').Value = '''1'
$ws.Range('E114').Value = '''True'
$ws.Range('F114').Value = '''sgd'
$ws.Range('G114').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H114').Value = '''0.01'
$ws.Range('I114').Value = '''0.9'
$ws.Range('J114').Value = '''<function exp_lr_scheduler at 0x7f0d4bb63840>'
$ws.Range('K114').Value = '''10'
$ws.Range('L114').Value = '''True'
$ws.Range('M114').Value = '''64'
$ws.Range('N114').Value = 13
$ws.Range('O114').Value = 0.0008734013459980413
$ws.Range('P114').Value = 0.0009759009182453155
$ws.Range('Q114').Value = 0.517490518523854
$ws.Range('R114').Value = 0.523
$ws.Range('S114').Value = 0.834661159632885
$ws.Range('T114').Value = 0.8689073598491384
$ws.Range('U114').Value = 13
$ws.Range('V114').Value = 0.517490518523854
$ws.Range('W114').Value = 0.523
$ws.Range('X114').Value = 0.834661159632885
$ws.Range('Y114').Value = 0.8689073598491384

# Row 115
$ws.Range('A115').Value = '''September28  18:12:40'
$ws.Range('B115').Value = '''resnet34'
$ws.Range('C115').Value = '''1000'
$ws.Range('D115').Value = '''1'
$ws.Range('E115').Value = '''True'
$ws.Range('F115').Value = '''sgd'
$ws.Range('G115').Value = '''0.0xsingle + 1.0Xmulti'
$ws.Range('H115').Value = '''0.01'
$ws.Range('I115').Value = '''0.9'
$ws.Range('J115').Value = '''<function exp_lr_scheduler at 0x7f4b85717ea0>'
$ws.Range('K115').Value = '''10'
$ws.Range('L115').Value = '''True'
$ws.Range('M115').Value = '''64'
$ws.Range('S115').Value = ''' '

# Row 116
$ws.Range('A116').Value = '''September28  18:13:49'
$ws.Range('B116').Value = '''resnet34'
$ws.Range('C116').Value = '''1000'
$ws.Range('D116').Value = '''1'
$ws.Range('E116').Value = '''True'
$ws.Range('F116').Value = '''sgd'
$ws.Range('G116').Value = '''0.0xsingle + 1.0Xmulti'
$ws.Range('H116').Value = '''0.01'
$ws.Range('I116').Value = '''0.9'
$ws.Range('J116').Value = '''<function exp_lr_scheduler at 0x7f4b85717ea0>'
$ws.Range('K116').Value = '''10'
$ws.Range('L116').Value = '''True'
$ws.Range('M116').Value = '''16'
$ws.Range('N116').Value = 5
$ws.Range('O116').Value = 0.01609985243218641
$ws.Range('P116').Value = 0.02417135745921033
$ws.Range('Q116').Value = 0.4557291666666667
$ws.Range('R116').Value = 0.4133738601823708
$ws.Range('S116').Value = 1.006489360765097
$ws.Range('T116').Value = 1.320863676184575
$ws.Range('U116').Value = 5
$ws.Range('V116').Value = 0.4557291666666667
$ws.Range('W116').Value = 0.4133738601823708
$ws.Range('X116').Value = 1.006489360765097
$ws.Range('Y116').Value = 1.320863676184575

# Row 117
$ws.Range('A117').Value = '''September28  18:21:34'
$ws.Range('B117').Value = '''resnet34'
$ws.Range('C117').Value = '''1000'
$ws.Range('D117').Value = '''1'
$ws.Range('E117').Value = '''True'
$ws.Range('F117').Value = '''sgd'
$ws.Range('G117').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H117').Value = '''0.01'
$ws.Range('I117').Value = '''0.9'
$ws.Range('J117').Value = '''<function exp_lr_scheduler at 0x7fb8d691f8c8>'
$ws.Range('K117').Value = '''10'
$ws.Range('L117').Value = '''True'
$ws.Range('M117').Value = '''64'
$ws.Range('S117').Value = ''' '

# Row 118
$ws.Range('A118').Value = '''September28  18:32:14'
$ws.Range('B118').Value = '''resnet34'
$ws.Range('C118').Value = '''1000'
$ws.Range('D118').Value = '''1'
$ws.Range('E118').Value = '''True'
$ws.Range('F118').Value = '''sgd'
$ws.Range('G118').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H118').Value = '''0.01'
$ws.Range('I118').Value = '''0.9'
$ws.Range('J118').Value = '''<function exp_lr_scheduler at 0x7fb8d691f8c8>'
$ws.Range('K118').Value = '''10'
$ws.Range('L118').Value = '''True'
$ws.Range('M118').Value = '''64'
$ws.Range('S118').Value = ''' '

# Row 119
$ws.Range('A119').Value = '''September28  18:46:01'
$ws.Range('B119').Value = '''resnet34'
$ws.Range('C119').Value = '''1000'
$ws.Range('D119').Value = '''1'
$ws.Range('E119').Value = '''True'
$ws.Range('F119').Value = '''sgd'
$ws.Range('G119').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H119').Value = '''0.01'
$ws.Range('I119').Value = '''0.9'
$ws.Range('J119').Value = '''<function exp_lr_scheduler at 0x7fde41d558c8>'
$ws.Range('K119').Value = '''10'
$ws.Range('L119').Value = '''True'
$ws.Range('M119').Value = '''64'
$ws.Range('S119').Value = ''' '

# Row 120
$ws.Range('A120').Value = '''September28  18:47:56'
$ws.Range('B120').Value = '''resnet34'
$ws.Range('C120').Value = '''1000'
$ws.Range('D120').Value = '''1'
$ws.Range('E120').Value = '''True'
$ws.Range('F120').Value = '''sgd'
$ws.Range('G120').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H120').Value = '''0.01'
$ws.Range('I120').Value = '''0.9'
$ws.Range('J120').Value = '''<function exp_lr_scheduler at 0x7fde41d558c8>'
$ws.Range('K120').Value = '''10'
$ws.Range('L120').Value = '''True'
$ws.Range('M120').Value = '''64'
$ws.Range('S120').Value = ''' '

# Row 121
$ws.Range('A121').Value = '''September28  18:49:10'
$ws.Range('B121').Value = '''resnet34'
$ws.Range('C121').Value = '''1000'
$ws.Range('D121').Value = '''1'
$ws.Range('E121').Value = '''True'
$ws.Range('F121').Value = '''sgd'
$ws.Range('G121').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H121').Value = '''0.01'
$ws.Range('I121').Value = '''0.9'
$ws.Range('J121').Value = '''<function exp_lr_scheduler at 0x7fde41d558c8>'
$ws.Range('K121').Value = '''10'
$ws.Range('L121').Value = '''True'
$ws.Range('M121').Value = '''64'
$ws.Range('S121').Value = ''' '

# Row 122
$ws.Range('A122').Value = '''October04  22:26:39'
$ws.Range('B122').Value = '''resnet34'
$ws.Range('C122').Value = '''1000'
$ws.Range('D122').Value = '''1'
$ws.Range('E122').Value = '''True'
$ws.Range('F122').Value = '''sgd'
$ws.Range('G122').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H122').Value = '''0.01'
$ws.Range('I122').Value = '''0.9'
$ws.Range('J122').Value = '''<function exp_lr_scheduler at 0x7fd3688c7a60>'
$ws.Range('K122').Value = '''10'
$ws.Range('L122').Value = '''True'
$ws.Range('M122').Value = '''64'
$ws.Range('S122').Value = ''' '

# Row 123
$ws.Range('A123').Value = '''October04  22:28:17'
$ws.Range('B123').Value = '''resnet34'
$ws.Range('C123').Value = '''1000'
$ws.Range('D123').Value = '''1'
$ws.Range('E123').Value = '''True'
$ws.Range('F123').Value = '''sgd'
$ws.Range('G123').Value = '''1.0xsingle + 0.0Xmulti'
$ws.Range('H123').Value = '''0.01'
$ws.Range('I123').Value = '''0.9'
$ws.Range('J123').Value = '''<function exp_lr_scheduler at 0x7f303bf42a60>'
$ws.Range('K123').Value = '''10'
$ws.Range('L123').Value = '''True'
$ws.Range('M123').Value = '''64'
$ws.Range('N123').Value = 6
$ws.Range('O123').Value = 0.00121326309654557
$ws.Range('P123').Value = 0.001203183844685555
$ws.Range('Q123').Value = 0.5053425936442943
$ws.Range('R123').Value = 0.503
$ws.Range('S123').Value = 0.9370445903807563
$ws.Range('T123').Value = 0.8933084573650918
$ws.Range('U123').Value = 6
$ws.Range('V123').Value = 0.5053425936442943
$ws.Range('W123').Value = 0.503
$ws.Range('X123').Value = 0.9370445903807563
$ws.Range('Y123').Value = 0.8933084573650918
